$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns at the front for Username / Password.
#    Old A:F (Patient ID..Contact Information) shifts to C:H.
# ---------------------------------------------------------------------------
$ws.Range("A:B").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("H1").Value = "Email Address"
$ws.Range("I1").Value = "Phone Number"

# Copy the existing header style onto the new header cells so they match the
# look of the rest of the header row (bold font, box border, alignment).
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Data rows: Username / Password columns
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "patient1"
$ws.Range("A3").Value = "patient2"
$ws.Range("A4").Value = "patient3"

$ws.Range("B2").Value = "password"
$ws.Range("B3").Value = "password"
$ws.Range("B4").Value = "password"

# ---------------------------------------------------------------------------
# 4. Update Name column (now column D) with the new patient names
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "John Doe"
$ws.Range("D3").Value = "Jane Smith"
$ws.Range("D4").Value = "Alice Johnson"

# ---------------------------------------------------------------------------
# 5. Re-write the Date of Birth column (E) as literal text so it keeps
#    looking like "1980-05-14" instead of being auto-converted to a date.
# ---------------------------------------------------------------------------
$ws.Range("E2:E4").NumberFormat = "@"
$ws.Range("E2").Value = "1980-05-14"
$ws.Range("E3").Value = "1975-11-22"
$ws.Range("E4").Value = "1990-07-08"
$ws.Range("E2:E4").Style = "Normal"

# ---------------------------------------------------------------------------
# 6. New Phone Number column (I) - stored as real numbers
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = 12345678
$ws.Range("I3").Value = 87654321
$ws.Range("I4").Value = 45678912

# ---------------------------------------------------------------------------
# 7. Give the Phone Number header (I1) its own left/right-only thin border
#    style (distinct from the full box border used by the rest of row 1).
# ---------------------------------------------------------------------------
$ws.Range("I1").Borders.Item(7).LineStyle = 1
$ws.Range("I1").Borders.Item(10).LineStyle = 1
$ws.Range("I1").Borders.Item(8).LineStyle = 0
$ws.Range("I1").Borders.Item(9).LineStyle = 0

# ---------------------------------------------------------------------------
# 8. Cosmetic touch-ups that mirror the authored workbook
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()
$hWidth = $ws.Columns.Item(8).ColumnWidth()
$ws.Columns.Item(9).ColumnWidth = $hWidth

$ws.Application.ActiveWindow.Zoom = 138
$ws.Range("E6").Select()

$ws.PageSetup.Orientation = 1

Write-Output "done"
